$d = $word.ActiveDocument

$replacements = @(
    @{old="36×79=2844"; new="70×31=2170"},
    @{old="36×86=3096"; new="48×49=2352"},
    @{old="83×26=2158"; new="31×52=1612"},
    @{old="30×71=2130"; new="91×83=7553"},
    @{old="78×77=6006"; new="32×40=1280"},
    @{old="46×80=3680"; new="50×89=4450"},
    @{old="72×62=4464"; new="62×57=3534"},
    @{old="79×19=1501"; new="25×35=875"},
    @{old="78×43=3354"; new="53×18=954"},
    @{old="41×36=1476"; new="21×31=651"},
    @{old="57×19=1083"; new="49×23=1127"},
    @{old="68×74=5032"; new="89×85=7565"},
    @{old="84×70=5880"; new="61×80=4880"},
    @{old="77×11=847"; new="13×36=468"},
    @{old="94×39=3666"; new="60×11=660"},
    @{old="28×98=2744"; new="95×30=2850"},
    @{old="24×34=816"; new="11×89=979"},
    @{old="99×28=2772"; new="30×34=1020"},
    @{old="85×73=6205"; new="25×83=2075"},
    @{old="60×90=5400"; new="22×82=1804"},
    @{old="34×48=1632"; new="90×57=5130"},
    @{old="66×80=5280"; new="19×37=703"},
    @{old="29×48=1392"; new="38×23=874"},
    @{old="42×41=1722"; new="12×34=408"},
    @{old="13×66=858"; new="57×75=4275"}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
